# Generate Report for Handoff
# Regenerates the localization-status report for a new source file
# (411b169f-c9f0-4c23-8bba-683a7d49dffb.md replaces 8ac53b9b-0631-4105-8595-710dacf8b319.md)
# and refreshes the handoff timestamps / target-file bookkeeping columns.

$wb = $excel.ActiveWorkbook

$oldBase = "8ac53b9b-0631-4105-8595-710dacf8b319"
$newBase = "411b169f-c9f0-4c23-8bba-683a7d49dffb"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A2").Value = "$newBase.md"

# B2 carries a hyperlink; update its display text in place so the existing
# relationship (and therefore its target URL) is preserved.
$bLink = $wsOverview.Range("B2").Hyperlinks.Item(1)
$bLink.TextToDisplay = "e2e\$newBase.md"

$wsOverview.Range("G2").Value = "2016-09-03 05:05:34"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$aLinkZh = $wsZhCn.Range("A2").Hyperlinks.Item(1)
$aLinkZh.TextToDisplay = "$newBase.md"

$wsZhCn.Range("G2").Value = "$newBase.97e178b1fb9a5db293eb37f30a3ae59a6d47818c.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 05:05:29"

# Latest Target File / Latest Handback File are no longer known for this
# handoff cycle - clear the values and drop the stale hyperlink on I2.
$wsZhCn.Range("I2").Hyperlinks.Item(1).Delete()
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("J2").Style = "Normal"

$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsZhCn.Columns.Item(10).ColumnWidth = 20.872143700009268

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$aLinkDe = $wsDeDe.Range("A2").Hyperlinks.Item(1)
$aLinkDe.TextToDisplay = "$newBase.md"

$wsDeDe.Range("G2").Value = "$newBase.97e178b1fb9a5db293eb37f30a3ae59a6d47818c.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-03 05:05:34"

$wsDeDe.Range("I2").Hyperlinks.Item(1).Delete()
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("J2").Style = "Normal"

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsDeDe.Columns.Item(10).ColumnWidth = 20.872143700009268
